# Update countries & provincias Spain
# - Refresh the COVID stats snapshot timestamp
# - Update case/death counters for the countries whose data changed
# - A handful of countries leap-frogged their neighbours in the
#   "Casos totales" ranking, so those rows swap places (name + stats)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- snapshot timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 2 de Julio de 2020 a las 11:44"

# --- per-row updates ------------------------------------------------------
# Each entry: row number, then a hashtable of column -> new value.
# Rows that changed rank also carry a new "A" (country name).
$updates = @(
    @{ Row = 7;   Cells = @{ B = 606907; C = 1687; D = 360378; E = 228669; G = 12; H = 17860 } },
    @{ Row = 21;  Cells = @{ B = 153277; C = 4019; D = 66442;  E = 84909;  G = 38; H = 1926 } },

    @{ Row = 31;  Cells = @{ A = "Indonesia"; B = 59394; C = 1624; D = 26667; E = 29740; G = 53; H = 2987 } },
    @{ Row = 32;  Cells = @{ A = "Ecuador";   B = 58257; D = 27887; E = 25794; H = 4576 } },

    @{ Row = 34;  Cells = @{ B = 50335 } },

    @{ Row = 40;  Cells = @{ A = "Oman";     B = 42555; C = 1361; D = 25318; E = 17049; G = 3; H = 188 } },
    @{ Row = 41;  Cells = @{ A = "Portugal"; B = 42454; D = 27798; E = 13077; H = 1579 } },

    @{ Row = 42;  Cells = @{ B = 38805; C = 294; D = 10673; E = 26858; G = 4; H = 1274 } },
    @{ Row = 43;  Cells = @{ B = 35146; C = 371; E = 11445; G = 15; H = 1492 } },
    @{ Row = 49;  Cells = @{ E = 5373; G = 1; H = 93 } },
    @{ Row = 60;  Cells = @{ B = 17941; C = 68; D = 16514; E = 722 } },

    @{ Row = 66;  Cells = @{ A = "Marruecos"; B = 12854; C = 218; D = 9052;  E = 3574; H = 228 } },
    @{ Row = 67;  Cells = @{ A = "Dinamarca"; B = 12794; D = 11693; E = 495;  H = 606 } },

    @{ Row = 74;  Cells = @{ B = 8643; C = 3; D = 8437; E = 85 } },
    @{ Row = 76;  Cells = @{ B = 7241; C = 5; E = 213 } },

    @{ Row = 103; Cells = @{ A = "Albania"; B = 2662; C = 82; D = 1559; E = 1034; G = 4; H = 69 } },
    @{ Row = 104; Cells = @{ A = "Mayotte"; B = 2643; D = 2341; E = 267; H = 35 } },

    @{ Row = 116; Cells = @{ B = 1825; C = 7; D = 1536; E = 211 } },

    @{ Row = 120; Cells = @{ A = "Eslovenia"; B = 1634; C = 21; D = 1384; E = 139; H = 111 } },
    @{ Row = 121; Cells = @{ A = "Zambia";    B = 1632; D = 1348; E = 254; H = 30 } },

    @{ Row = 133; Cells = @{ D = 947; E = 61 } },

    @{ Row = 163; Cells = @{ A = "Namibia"; C = 8; D = 24;  E = 269; H = 0 } },
    @{ Row = 164; Cells = @{ A = "Siria";   B = 293; D = 110; E = 174; H = 9 } },
    @{ Row = 165; Cells = @{ A = "Angola";  B = 291; D = 97;  E = 179; H = 15 } },

    @{ Row = 203; Cells = @{ A = "Santa Lucia" } },
    @{ Row = 204; Cells = @{ A = "Laos" } }
)

foreach ($u in $updates) {
    $row = $u.Row
    foreach ($col in $u.Cells.Keys) {
        $ws.Range("$col$row").Value = $u.Cells[$col]
    }
}
